$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GDPbES")

# Insert a new blank row at row 3 (pushes "natural gas nonpeaker" and everything
# below it down by one row, and Excel auto-adjusts any relative formula
# references that point at rows >= 3).
$ws.Rows.Item(3).Insert()

# Row 3 becomes "natural gas steam turbine" (new row), value 0 for 2015,
# filled right through column AK with formula "=$B3" like the other rows.
$ws.Range("A3").Value = "natural gas steam turbine"
$ws.Range("B3").Value = 0
$ws.Range("C3:AK3").Formula = "=`$B3"

# Row 4 was "natural gas nonpeaker" (now shifted down); rename it to
# "natural gas combined cycle" and set its value to 0 across the row.
$ws.Range("A4").Value = "natural gas combined cycle"
$ws.Range("B4").Value = 0
$ws.Range("C4:AK4").Formula = "=`$B4"

# Row 15 ("offshore wind") previously referenced B11 via formula; the diff
# shows it becomes a plain literal value (no formula) after the edit.
$ws.Range("B15:AK15").Value = 0

# Column A needs to be wide enough to fit the new, longer labels.
$ws.Columns.Item(1).ColumnWidth = 26.6
